$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "Dallas"
$ws.Cells.Item(2, 2).Value = 6
$ws.Cells.Item(2, 3).Value = 2325
$ws.Cells.Item(2, 4).Value = 387.5
$ws.Cells.Item(2, 5).Value = 1658
$ws.Cells.Item(2, 6).Value = 276.3
$ws.Cells.Item(2, 7).Value = 703
$ws.Cells.Item(2, 8).Value = 117.2
$ws.Cells.Item(2, 9).Value = 178
$ws.Cells.Item(2, 10).Value = 29.7

$ws.Cells.Item(3, 1).Value = "Atlanta"
$ws.Cells.Item(3, 2).Value = 5
$ws.Cells.Item(3, 3).Value = 1894
$ws.Cells.Item(3, 4).Value = 378.8
$ws.Cells.Item(3, 5).Value = 1197
$ws.Cells.Item(3, 6).Value = 239.4
$ws.Cells.Item(3, 7).Value = 756
$ws.Cells.Item(3, 8).Value = 151.19999999999999
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(3, 10).Value = 20

$ws.Cells.Item(4, 1).Value = "Buffalo"
$ws.Cells.Item(4, 2).Value = 6
$ws.Cells.Item(4, 3).Value = 2270
$ws.Cells.Item(4, 4).Value = 378.3
$ws.Cells.Item(4, 5).Value = 1429
$ws.Cells.Item(4, 6).Value = 238.2
$ws.Cells.Item(4, 7).Value = 906
$ws.Cells.Item(4, 8).Value = 151
$ws.Cells.Item(4, 9).Value = 167
$ws.Cells.Item(4, 10).Value = 27.8

$ws.Cells.Item(5, 1).Value = "Indianapolis"
$ws.Cells.Item(5, 2).Value = 6
$ws.Cells.Item(5, 3).Value = 2261
$ws.Cells.Item(5, 4).Value = 376.8
$ws.Cells.Item(5, 5).Value = 1511
$ws.Cells.Item(5, 6).Value = 251.8
$ws.Cells.Item(5, 7).Value = 791
$ws.Cells.Item(5, 8).Value = 131.80000000000001
$ws.Cells.Item(5, 9).Value = 194
$ws.Cells.Item(5, 10).Value = 32.299999999999997

$ws.Cells.Item(6, 1).Value = "L.A. Rams"
$ws.Cells.Item(6, 2).Value = 6
$ws.Cells.Item(6, 3).Value = 2250
$ws.Cells.Item(6, 4).Value = 375
$ws.Cells.Item(6, 5).Value = 1684
$ws.Cells.Item(6, 6).Value = 280.7
$ws.Cells.Item(6, 7).Value = 635
$ws.Cells.Item(6, 8).Value = 105.8
$ws.Cells.Item(6, 9).Value = 140
$ws.Cells.Item(6, 10).Value = 23.3

$ws.Cells.Item(7, 1).Value = "San Francisco"
$ws.Cells.Item(7, 2).Value = 6
$ws.Cells.Item(7, 3).Value = 2242
$ws.Cells.Item(7, 4).Value = 373.7
$ws.Cells.Item(7, 5).Value = 1838
$ws.Cells.Item(7, 6).Value = 306.3
$ws.Cells.Item(7, 7).Value = 493
$ws.Cells.Item(7, 8).Value = 82.2
$ws.Cells.Item(7, 9).Value = 125
$ws.Cells.Item(7, 10).Value = 20.8

$ws.Cells.Item(14, 1).Value = "Washington"
$ws.Cells.Item(14, 2).Value = 6
$ws.Cells.Item(14, 3).Value = 2074
$ws.Cells.Item(14, 4).Value = 345.7
$ws.Cells.Item(14, 5).Value = 1238
$ws.Cells.Item(14, 6).Value = 206.3
$ws.Cells.Item(14, 7).Value = 906
$ws.Cells.Item(14, 8).Value = 151
$ws.Cells.Item(14, 9).Value = 158
$ws.Cells.Item(14, 10).Value = 26.3

$ws.Cells.Item(15, 1).Value = "Chicago"
$ws.Cells.Item(15, 2).Value = 5
$ws.Cells.Item(15, 3).Value = 1693
$ws.Cells.Item(15, 4).Value = 338.6
$ws.Cells.Item(15, 5).Value = 1206
$ws.Cells.Item(15, 6).Value = 241.2
$ws.Cells.Item(15, 7).Value = 554
$ws.Cells.Item(15, 8).Value = 110.8
$ws.Cells.Item(15, 9).Value = 126
$ws.Cells.Item(15, 10).Value = 25.2

$ws.Cells.Item(16, 1).Value = "Denver"
$ws.Cells.Item(16, 2).Value = 6
$ws.Cells.Item(16, 3).Value = 2022
$ws.Cells.Item(16, 4).Value = 337
$ws.Cells.Item(16, 5).Value = 1277
$ws.Cells.Item(16, 6).Value = 212.8
$ws.Cells.Item(16, 7).Value = 781
$ws.Cells.Item(16, 8).Value = 130.19999999999999
$ws.Cells.Item(16, 9).Value = 130
$ws.Cells.Item(16, 10).Value = 21.7

$ws.Cells.Item(17, 1).Value = "New England"
$ws.Cells.Item(17, 2).Value = 6
$ws.Cells.Item(17, 3).Value = 2016
$ws.Cells.Item(17, 4).Value = 336
$ws.Cells.Item(17, 5).Value = 1522
$ws.Cells.Item(17, 6).Value = 253.7
$ws.Cells.Item(17, 7).Value = 549
$ws.Cells.Item(17, 8).Value = 91.5
$ws.Cells.Item(17, 9).Value = 150
$ws.Cells.Item(17, 10).Value = 25

$ws.Cells.Item(18, 1).Value = "Jacksonville"
$ws.Cells.Item(18, 2).Value = 6
$ws.Cells.Item(18, 3).Value = 1986
$ws.Cells.Item(18, 4).Value = 331
$ws.Cells.Item(18, 5).Value = 1324
$ws.Cells.Item(18, 6).Value = 220.7
$ws.Cells.Item(18, 7).Value = 744
$ws.Cells.Item(18, 8).Value = 124
$ws.Cells.Item(18, 9).Value = 139
$ws.Cells.Item(18, 10).Value = 23.2

$ws.Cells.Item(19, 1).Value = "Carolina"
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(19, 3).Value = 1985
$ws.Cells.Item(19, 4).Value = 330.8
$ws.Cells.Item(19, 5).Value = 1208
$ws.Cells.Item(19, 6).Value = 201.3
$ws.Cells.Item(19, 7).Value = 856
$ws.Cells.Item(19, 8).Value = 142.69999999999999
$ws.Cells.Item(19, 9).Value = 132
$ws.Cells.Item(19, 10).Value = 22

$ws.Cells.Item(20, 1).Value = "N.Y. Giants"
$ws.Cells.Item(20, 2).Value = 6
$ws.Cells.Item(20, 3).Value = 1969
$ws.Cells.Item(20, 4).Value = 328.2
$ws.Cells.Item(20, 5).Value = 1294
$ws.Cells.Item(20, 6).Value = 215.7
$ws.Cells.Item(20, 7).Value = 757
$ws.Cells.Item(20, 8).Value = 126.2
$ws.Cells.Item(20, 9).Value = 121
$ws.Cells.Item(20, 10).Value = 20.2
